$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the header row (row 1) with new English labels (mockup school data headers)
$ws.Range("A1").Value = "index"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "name-surname"
$ws.Range("D1").Value = "class-room"
$ws.Range("E1").Value = "father-name"
$ws.Range("F1").Value = "father-phone"
$ws.Range("G1").Value = "father-email"
$ws.Range("H1").Value = "mother-name"
$ws.Range("I1").Value = "mother-phone"
$ws.Range("J1").Value = "mother-email"
$ws.Range("K1").Value = "address"

# Update the active selection to C4 (matches sheetView selection change)
$ws.Range("C4").Select()

$wb.Save()
